$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.401.12'
$ws.Range("E2").Value = '  +1.10%  '
$ws.Range("D3").Value = '2.913.44'
$ws.Range("E3").Value = '  +3.70%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '352.71'
$ws.Range("E5").Value = '  -0.76%  '
$ws.Range("D6").Value = '112.25'
$ws.Range("E6").Value = '  +0.72%  '
$ws.Range("D7").Value = '0.558'
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '0.634'
$ws.Range("E9").Value = '  +0.46%  '
$ws.Range("D10").Value = '39.92'
$ws.Range("E10").Value = '  -1.08%  '
$ws.Range("E11").Value = '  +3.05%  '
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").Value = '19.88'
$ws.Range("E13").Value = '  -0.67%  '
$ws.Range("E14").Value = '  +0.38%  '
$ws.Range("D15").Value = '3.371.52'
$ws.Range("E15").Value = '  +3.65%  '
$ws.Range("E16").Value = '  +6.83%  '
$ws.Range("D17").Value = '2.915.37'
$ws.Range("E17").Value = '  +3.74%  '
$ws.Range("D18").Value = '52.414.79'
$ws.Range("E18").Value = '  +1.12%  '
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("E20").Value = '  +3.19%  '
$ws.Range("D21").Value = '14.16'
$ws.Range("E21").Value = '  +3.67%  '
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").Value = '70.96'
$ws.Range("E23").Value = '  +0.70%  '
$ws.Range("D24").Value = '270.24'
$ws.Range("E24").Value = '  +0.48%  '
$ws.Range("D25").Value = '2.77'
$ws.Range("E25").Value = '  -0.17%  '
$ws.Range("E26").Value = '  +2.20%  '
$ws.Range("D27").Value = '0.169'
$ws.Range("E27").Value = '  +4.23%  '
$ws.Range("E28").Value = '  -0.19%  '
$ws.Range("D29").Value = '10.66'
$ws.Range("E29").Value = '  +2.55%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").Value = '6.38'
$ws.Range("E30").Value = '  +12.90%  '
$ws.Range("B31").Value = 'InjectiveProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D31").Value = '37.86'
$ws.Range("E31").Value = '  -2.18%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.60'
$ws.Range("E32").Value = '  +7.70%  '
$ws.Range("B33").Value = 'Toncoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D33").Value = '2.25'
$ws.Range("E33").Value = '  +0.48%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.0982'
$ws.Range("E34").Value = '  +11.31%  '
$ws.Range("E35").Value = '  +1.89%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0450'
$ws.Range("E36").Value = '  +1.29%  '
$ws.Range("E37").Value = '  -0.05%  '
$ws.Range("D38").Value = '3.32'
$ws.Range("E38").Value = '  +5.45%  '
$ws.Range("E39").Value = '  -0.26%  '
$ws.Range("D40").Value = '2.07'
$ws.Range("E40").Value = '  +2.46%  '
$ws.Range("D41").Value = '2.85'
$ws.Range("E41").Value = '  +13.59%  '
$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").Value = '0.117'
$ws.Range("E42").Value = '  +1.40%  '
$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '23.51'
$ws.Range("E43").Value = '  +6.04%  '
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").Value = '120.82'
$ws.Range("E44").Value = '  +0.22%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.60'
$ws.Range("E45").Value = '  +7.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.20'
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("D47").Value = '3.54'
$ws.Range("E47").Value = '  +3.44%  '
$ws.Range("D48").Value = '2.200.03'
$ws.Range("E48").Value = '  +4.31%  '
$ws.Range("E49").Value = '  +24.30%  '
$ws.Range("D50").Value = '0.0343'
$ws.Range("E50").Value = '  +12.74%  '
$ws.Range("D51").Value = '0.972'
